# Update regression coefficient text values on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B3: -2.821*** -> -2.82***
$ws.Range("B3").Value = "-2.82***"

# C2: -0.012* -> -0.01*
$ws.Range("C2").Value = "-0.01*"

# C3: -0.467*** -> -0.47***
$ws.Range("C3").Value = "-0.47***"
